$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct revised values for the last three existing rows (314-316) ---
$ws.Range("C314:F314").Value = 2078216858000
$ws.Range("C315:F315").Value = 2082183969000
$ws.Range("C316:F316").Value = 2118202312000

# --- Append new rows 317-319 ---
$ws.Range("A317").Value = 44986.45833333334
$ws.Range("B317").Value = "ECONOMICS:PLM2"
$ws.Range("C317:F317").Value = 2121975670000
$ws.Range("G317").Value = 0

$ws.Range("A318").Value = 45017.45833333334
$ws.Range("B318").Value = "ECONOMICS:PLM2"
$ws.Range("C318:F318").Value = 2135028350000
$ws.Range("G318").Value = 0

$ws.Range("A319").Value = 45047.41666666666
$ws.Range("B319").Value = "ECONOMICS:PLM2"
$ws.Range("C319:F319").Value = 2140971740000
$ws.Range("G319").Value = 0

# Copy the date/time number formatting (style) from row 316's A cell onto
# the newly added A317:A319 cells, matching the existing column A styling.
$ws.Range("A316").Copy()
$ws.Range("A317:A319").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the worksheet dimension to reflect the new extent of used cells.
$ws.Range("A1:G319").Select()
